$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("git-manual")

$ws.Range("A17").Value = "git diff --cached HEAD~"
$ws.Range("B17").Value = "比较index与HEAD~ commit的差异"
$ws.Range("A18").Value = "git diff HEAD~ HEAD"
$ws.Range("B18").Value = "输出HEAD较之于HEAD~ commit的差异"

$ws.Range("B18").Select()
